$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)   # "总计"
$q1Sheet    = $wb.Worksheets.Item(2)   # "2022-Q1" (will stay, just shifts position)
$q4_21Sheet = $wb.Worksheets.Item(3)   # "2021-Q4" (will stay, just shifts position)

# ---------------------------------------------------------------------------
# Step 1: Update "总计" (totals) sheet - a new 2022-Q4 row is inserted as the
# first data row, pushing the existing 2022-Q1 / 2021-Q4 rows down by one.
# ---------------------------------------------------------------------------

# Copy the formatting (style index) of the existing last data row's A-cell
# down to the new row 4 before we start overwriting values, so the
# sequential-index column A keeps the same cell style used by the other rows.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)

# Shift the old row 3 (2021-Q4) down into row 4.
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q4"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0

# Shift the old row 2 (2022-Q1) down into row 3.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q1"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.05

# Write the brand new 2022-Q4 totals into row 2.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.1

# ---------------------------------------------------------------------------
# Step 2: Create the brand-new "2022-Q4" holdings worksheet, positioned right
# after "总计" (so the sheet order becomes 总计, 2022-Q4, 2022-Q1, 2021-Q4).
#
# Rather than Worksheets.Add() (which produces a bare sheet lacking the
# sheetPr/pageMargins seen on every other sheet in this workbook), duplicate
# the existing "2022-Q1" sheet - it already carries the right sheetPr,
# pageMargins, column layout and header text/styles - then overwrite its
# data rows with the 2022-Q4 numbers.
# ---------------------------------------------------------------------------
$q1Sheet.Copy($null, $totalSheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# Extend the index-column style (s="2") from row 2 down into row 3.
$newSheet.Range("A2").Copy()
$newSheet.Range("A3").PasteSpecial(-4122)

# The source sheet stores fund code / size / position figures as literal
# TEXT (inlineStr), not numbers - e.g. "005457" must keep its leading zero.
# Plain `.Value = "005457"` gets auto-coerced to the number 5457 by this
# engine, so route the numeric-looking strings through helper cells holding
# `="literal"` formulas (always typed as text) and paste-special the
# resulting values in, which preserves the text cell type.
$newSheet.Range("Z1").Formula  = '="005457"'
$newSheet.Range("Z2").Formula  = '="5.08"'
$newSheet.Range("Z3").Formula  = '="94.36"'
$newSheet.Range("Z4").Formula  = '="1.47"'
$newSheet.Range("Z5").Formula  = '="0.0747"'
$newSheet.Range("Z6").Formula  = '="008851"'
$newSheet.Range("Z7").Formula  = '="2.37"'
$newSheet.Range("Z8").Formula  = '="71.22"'
$newSheet.Range("Z9").Formula  = '="1.21"'
$newSheet.Range("Z10").Formula = '="0.0287"'

$newSheet.Range("Z1").Copy()
$newSheet.Range("B2").PasteSpecial(-4163)
$newSheet.Range("Z2").Copy()
$newSheet.Range("D2").PasteSpecial(-4163)
$newSheet.Range("Z3").Copy()
$newSheet.Range("E2").PasteSpecial(-4163)
$newSheet.Range("Z4").Copy()
$newSheet.Range("F2").PasteSpecial(-4163)
$newSheet.Range("Z5").Copy()
$newSheet.Range("G2").PasteSpecial(-4163)
$newSheet.Range("Z6").Copy()
$newSheet.Range("B3").PasteSpecial(-4163)
$newSheet.Range("Z7").Copy()
$newSheet.Range("D3").PasteSpecial(-4163)
$newSheet.Range("Z8").Copy()
$newSheet.Range("E3").PasteSpecial(-4163)
$newSheet.Range("Z9").Copy()
$newSheet.Range("F3").PasteSpecial(-4163)
$newSheet.Range("Z10").Copy()
$newSheet.Range("G3").PasteSpecial(-4163)

# Clean up the helper column.
$newSheet.Range("Z1:Z10").ClearContents()

# Fill in the remaining plain-text / numeric cells.
$newSheet.Range("A2").Value = 0
$newSheet.Range("C2").Value = "景顺长城量化小盘股票"
$newSheet.Range("H2").Value = 4

$newSheet.Range("A3").Value = 1
$newSheet.Range("C3").Value = "景顺长城量化对冲策略三个月定期开放灵活配置混合"
$newSheet.Range("H3").Value = 3

# ---------------------------------------------------------------------------
# Step 3: Restore the originally-active "2021-Q4" sheet (now the 4th tab) as
# the selected one, since inserting/copying sheets changes the active tab.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q4").Activate()
